$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sheet title / "through" date from Oct 04 to Oct 05 ---
$ws.Name = "Through 2021-10-05"
$ws.Range("B1").Value = "October 2021 (through October 05)"

# --- Update per-neighborhood carjacking counts with the new day's data ---

# Garfield Park (row 2): Oct 2021 +1, Oct 2020 +1
$ws.Range("B2").Value = 3
$ws.Range("L2").Value = 4

# Austin (row 3): Oct 2018 +1, Oct 2017 +1
$ws.Range("AF3").Value = 3
$ws.Range("AP3").Value = 2

# Wicker Park (row 7): new Oct 2020 count
$ws.Range("L7").Value = 1

# Roseland (row 10): Oct 2020 +1, new Oct 2018 count
$ws.Range("L10").Value = 2
$ws.Range("AF10").Value = 1

# South Shore (row 13): new Oct 2020 count
$ws.Range("L13").Value = 1

# Ashburn (row 21): new Oct 2019 count
$ws.Range("V21").Value = 1

# Grand Boulevard (row 22): Jan 2021 +1
$ws.Range("K22").Value = 9

# River North (row 25): new Oct 2016 count
$ws.Range("AZ25").Value = 1

# Lincoln Park (row 31): new Oct 2020 count
$ws.Range("L31").Value = 1

# Hyde Park (row 39): new Oct 2015 count
$ws.Range("BJ39").Value = 1

# Washington Heights (row 48): new Oct 2021 count
$ws.Range("B48").Value = 1

# Albany Park (row 59): new Oct 2020 count
$ws.Range("L59").Value = 1

# Calumet Heights (row 65): new Oct 2020 count
$ws.Range("L65").Value = 1

# East Village (row 68): new Oct 2021 count
$ws.Range("B68").Value = 1

# Magnificent Mile (row 80): new Oct 2021 count
$ws.Range("B80").Value = 1
